$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3962093333333334
$ws.Range("H2").Value = 1.188628
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.618716333333334
$ws.Range("N2").Value = 7.856149000000001
$ws.Range("O2").Value = 0.07115908183301342
$ws.Range("P2").Value = 0.07115908183301341
$ws.Range("Q2").Value = 1.037559852619111
$ws.Range("R2").Value = 9.338038673572001
$ws.Range("S2").Value = 0.07115908183301342
$ws.Range("T2").Value = 0.07115908183301341

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3962093333333334
$ws.Range("H3").Value = 1.188628
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.503283
$ws.Range("N3").Value = 46.509849
$ws.Range("O3").Value = 0.4212748702999519
$ws.Range("P3").Value = 0.4212748702999519
$ws.Range("Q3").Value = 6.142545421908001
$ws.Range("R3").Value = 55.28290879717201
$ws.Range("S3").Value = 0.4212748702999519
$ws.Range("T3").Value = 0.4212748702999519

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3962093333333334
$ws.Range("H4").Value = 1.188628
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.67887366666666
$ws.Range("N4").Value = 56.036621
$ws.Range("O4").Value = 0.5075660478670347
$ws.Range("P4").Value = 0.5075660478670347
$ws.Range("Q4").Value = 7.400744082887555
$ws.Range("R4").Value = 66.606696745988
$ws.Range("S4").Value = 0.5075660478670347
$ws.Range("T4").Value = 0.5075660478670347
